$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    $txt = ""
    if ($sh.HasTextFrame) {
        $txt = $sh.TextFrame.TextRange.Text
    }
    Write-Host "$i : id=$($sh.Id) name=$($sh.Name) type=$($sh.Type) L=$($sh.Left) T=$($sh.Top) W=$($sh.Width) H=$($sh.Height) text=[$txt]"
}
